# Change schottky diode on fixed,10V,24mA,275nm design.
# Update placement coordinates (Mid X / Mid Y) for parts C2, D2, R5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 -> Designator C2
$ws.Range("B3").Value = 65.4902
$ws.Range("C3").Value = -52.8398

# Row 7 -> Designator D2
$ws.Range("B7").Value = 65.775
$ws.Range("C7").Value = -54.15

# Row 16 -> Designator R5
$ws.Range("B16").Value = 65.4602
$ws.Range("C16").Value = -51.773

# Update the selection state to match the recorded edit (B3:C3, active cell B3)
$ws.Range("B3:C3").Select()
